# Fig33.xlsx - 2017-02-13 snapshot (chunk 30)
# Updates the "Short-Term Energy Outlook" title/source strings from
# January 2017 to February 2017, refreshes the trailing ~month of
# index/change data (rows 69-99) with revised STEO figures, and moves
# the forecast-divider marker from month 48.5 to 49.5 (one more month
# of "actual" history rolled in).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fig33")

# --- Title (A2) and chart source note (A100) ------------------------------
$ws.Range("A2").Value   = "Short-Term Energy Outlook, February 2017"
$ws.Range("A100").Value = "Source: Short-Term Energy Outlook, February 2017."

# --- Revised monthly data, rows 69-99 (column B = Index, column C = Change) -
$newB = @(104.2222, 104.54170000000001, 104.4537, 104.24290000000001, 104.3185, 103.85890000000001, 104.37854197999999, 104.72867037, 104.9034, 105.0497, 105.06010000000001, 105.2303, 105.45269999999999, 105.8045, 106.0735, 106.3368, 106.568, 106.84, 107.12609999999999, 107.4542, 107.748, 108.0352, 108.2919, 108.5839, 108.8873, 109.2546, 109.54130000000001, 109.8001, 110.012, 110.22880000000001, 110.4315)

$newC = @(-0.0060814477221999997, -0.0088532408000000003, -0.010651810078, -0.010106621389999999, -0.0080483126974999989, -0.0060122254326000005, 0.0032038188721000001, 0.0017137372284999999, 0.0046867209999999996, 0.015703769999999999, 0.01176464, 0.01481854, 0.0118064, 0.012079319999999999, 0.015507169999999999, 0.020086989999999999, 0.021564119999999999, 0.028702969999999998, 0.026323050000000001, 0.02602471, 0.027116699999999997, 0.028420209999999998, 0.030761189999999997, 0.031869330000000001, 0.032570130000000003, 0.032608290000000005, 0.032692990000000005, 0.032568559999999996, 0.032317100000000001, 0.031718400000000001, 0.030855570000000002)

$startRow = 69
for ($i = 0; $i -lt $newB.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $newB[$i]
    $ws.Cells.Item($r, 3).Value = $newC[$i]
}

# --- Forecast-divider marker (scatter series x-values), rows 103-104 ------
$ws.Range("A103").Value = 49.5
$ws.Range("A104").Value = 49.5

$wb.Save()
